$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.544.33"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.666.75"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.01"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4797"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2631"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06171"
$ws.Range("E9").Value = "  +2.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07086"
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.664.66"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.81"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5881"
$ws.Range("E13").Value = "  -5.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.370"
$ws.Range("E14").Value = "  -4.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.07"
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.0000"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.0000"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.535.80"
$ws.Range("E18").Value = "  +2.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006746"
$ws.Range("E19").Value = "  +2.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.877.04"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.414"
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.731"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.276"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.86"
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.02"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.389"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "104.90"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.713"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.975"
$ws.Range("E30").Value = "  +5.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07755"
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.642"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9989"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04223"
$ws.Range("E34").Value = "  -7.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.600"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6112"
$ws.Range("E36").Value = "  +5.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9490"
$ws.Range("E37").Value = "  +1.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.592"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8634"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9996"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.849"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01466"
$ws.Range("E42").Value = "  -6.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.16"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3765"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.860"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1120"
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.205"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.375"
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("E51").Value = "  -0.03%  "
